$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab name + workbook sheet entry) from "Through 2021-09-10" to "Through 2021-09-11"
$ws.Name = "Through 2021-09-11"

# Update the "September (through 09-10)" label to "September (through 09-11)"
$ws.Range("A10").Value = "September (through 09-11)"

# Update September row (row 10) values for each year column
$ws.Range("B10").Value = 10
$ws.Range("C10").Value = 18
$ws.Range("D10").Value = 27
$ws.Range("E10").Value = 20
$ws.Range("F10").Value = 25
$ws.Range("G10").Value = 37
$ws.Range("H10").Value = 53

# Update Total row (row 11) values for each year column
$ws.Range("B11").Value = 204
$ws.Range("C11").Value = 399
$ws.Range("D11").Value = 578
$ws.Range("E11").Value = 510
$ws.Range("F11").Value = 374
$ws.Range("G11").Value = 821
$ws.Range("H11").Value = 1124
